$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fn1"
$ws.Range("C2").Value = "Mag"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.116717
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.09045733333333333
$ws.Range("N2").Value = 0.271372
$ws.Range("O2").Value = 0.0574491187297735
$ws.Range("P2").Value = 0.0574491187297735
$ws.Range("Q2").Value = 2.445867302858222
$ws.Range("R2").Value = 22.012805725724
$ws.Range("S2").Value = 0.004076697595240886
$ws.Range("T2").Value = 0.004076697595240887

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fn1"
$ws.Range("C3").Value = "Mag"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.116717
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.484107
$ws.Range("N3").Value = 4.452321
$ws.Range("O3").Value = 0.9425508812702265
$ws.Range("P3").Value = 0.9425508812702265
$ws.Range("Q3").Value = 40.12862917223966
$ws.Range("R3").Value = 361.157662550157
$ws.Range("S3").Value = 0.06688518459509639
$ws.Range("T3").Value = 0.0668851845950964

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fn1"
$ws.Range("C4").Value = "Mag"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 345.566579
$ws.Range("H4").Value = 1036.699737
$ws.Range("I4").Value = 0.9069174311350353
$ws.Range("J4").Value = 0.9069174311350354
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.09045733333333333
$ws.Range("N4").Value = 0.271372
$ws.Range("O4").Value = 0.0574491187297735
$ws.Range("P4").Value = 0.0574491187297735
$ws.Range("Q4").Value = 31.25903122546267
$ws.Range("R4").Value = 281.331281029164
$ws.Range("S4").Value = 0.05210160717937783
$ws.Range("T4").Value = 0.05210160717937783

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fn1"
$ws.Range("C5").Value = "Mag"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.484107
$ws.Range("N5").Value = 4.452321
$ws.Range("O5").Value = 0.9425508812702265
$ws.Range("P5").Value = 0.9425508812702265
$ws.Range("Q5").Value = 512.857778859953
$ws.Range("R5").Value = 4615.720009739577
$ws.Range("S5").Value = 0.8548158239556575
$ws.Range("T5").Value = 0.8548158239556576

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Fn1"
$ws.Range("C6").Value = "Mag"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.428738666666668
$ws.Range("H6").Value = 25.286216
$ws.Range("I6").Value = 0.0221206866746274
$ws.Range("J6").Value = 0.02212068667462741
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.09045733333333333
$ws.Range("N6").Value = 0.271372
$ws.Range("O6").Value = 0.0574491187297735
$ws.Range("P6").Value = 0.0574491187297735
$ws.Range("Q6").Value = 0.7624412231502223
$ws.Range("R6").Value = 6.861971008352
$ws.Range("S6").Value = 0.001270813955154788
$ws.Range("T6").Value = 0.001270813955154788

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Fn1"
$ws.Range("C7").Value = "Mag"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.428738666666668
$ws.Range("H7").Value = 25.286216
$ws.Range("I7").Value = 0.0221206866746274
$ws.Range("J7").Value = 0.02212068667462741
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.484107
$ws.Range("N7").Value = 4.452321
$ws.Range("O7").Value = 0.9425508812702265
$ws.Range("P7").Value = 0.9425508812702265
$ws.Range("Q7").Value = 12.50915005637067
$ws.Range("R7").Value = 112.582350507336
$ws.Range("S7").Value = 0.02084987271947261
$ws.Range("T7").Value = 0.02084987271947262
